$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case row 14: highlighted command cell (same style as B5/B11)
$ws.Range("B14").Value = "cat split_* | grep 'block->' >test1 && cat test1"
$ws.Range("B14").Interior.Color = 65535

# New test case row 16: highlighted command cell (same style as B5/B11)
$ws.Range("B16").Value = 'echo "hello""bye"'
$ws.Range("B16").Interior.Color = 65535

# Page setup becomes explicit (paper size 9 = A4, portrait orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active selection to reflect where the user ended up
[void]$ws.Range("M16").Select()
